$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing column C, shifting the old
# C:G (URL/IT-consulting detail) columns over to D:H. Excel fills the new
# column's cell styles in from the neighbouring column automatically.
$ws.Columns("C").Insert()

# Header for the new column.
$ws.Range("C3").Value2 = "URL Website"

# URL values for each of the 14 data rows (one per IT consulting company).
$ws.Range("C4").Value2 = "https://www.onetrust.com/"
$ws.Range("C5").Value2 = "https://brightbridgesolutions.com/"
$ws.Range("C6").Value2 = "https://www.kellton.com/"
$ws.Range("C7").Value2 = "https://mirus-it.co.uk/"
$ws.Range("C8").Value2 = "https://www.russellreynolds.com/en/"
$ws.Range("C9").Value2 = "https://www.bain.com/offices/london/"
$ws.Range("C10").Value2 = "https://www.itransition.com/"
$ws.Range("C11").Value2 = "https://www.webbytech.co.uk/"
$ws.Range("C12").Value2 = "https://trilan-it.com/"
$ws.Range("C13").Value2 = "https://www.lumen.com/en-us/home.html"
$ws.Range("C14").Value2 = "https://conosco.com/"
$ws.Range("C15").Value2 = "https://www.transparity.com/latest-news/transparity-acquire-microsoft-azure-development-specialist-ballard-chalmers/"
$ws.Range("C16").Value2 = "https://www.ecl.co.uk/"
$ws.Range("C17").Value2 = "https://www.topdesk.com/en/"

# The first URL (C4) was turned into a real hyperlink at some point, which
# is what introduced the "Hyperlink" named cell style / underlined themed
# font into the workbook's style table. Recreate that, then remove the
# hyperlink relationship itself (keeping only the visual formatting it left
# behind on the cell) so the final workbook matches the target - styled
# like a link, but without a live hyperlink.
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.onetrust.com/")
$ws.Hyperlinks.Delete()

# Leave the selection where the edit finished.
$ws.Range("C17").Select()
